$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.325.52"
$ws.Range("E2").Value = "  +2.77%  "

# Row 3
$ws.Range("D3").Value = "2.094.05"
$ws.Range("E3").Value = "  +3.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "

# Row 6
$ws.Range("E6").Value = "  -0.35%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +14.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.372"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.04%  "

# Row 11
$ws.Range("E11").Value = "  +3.89%  "

# Row 12
$ws.Range("E12").Value = "  +7.31%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.38%  "

# Row 14
$ws.Range("D14").Value = "2.399.71"
$ws.Range("E14").Value = "  +4.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.835"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.82%  "

# Row 16
$ws.Range("D16").Value = "2.092.98"
$ws.Range("E16").Value = "  +4.01%  "

# Row 17
$ws.Range("E17").Value = "  +4.40%  "

# Row 18
$ws.Range("D18").Value = "37.220.73"
$ws.Range("E18").Value = "  +3.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.86%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  +2.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.13%  "

# Row 23
$ws.Range("E23").Value = "  +6.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.94%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.63%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.30%  "

# Row 30
$ws.Range("E30").Value = "  +0.65%  "

# Row 31
$ws.Range("E31").Value = "  +27.36%  "

# Row 32
$ws.Range("E32").Value = "  +2.56%  "

# Row 33
$ws.Range("B33").Value = "Gas"
$ws.Range("C33").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0611"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.57%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0914"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.97%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("E37").Value = "  +8.86%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "

# Row 40
$ws.Range("E40").Value = "  -1.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.95%  "

# Row 42
$ws.Range("E42").Value = "  +4.18%  "

# Row 43
$ws.Range("E43").Value = "  +6.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.93%  "

# Row 48
$ws.Range("D48").Value = "1.318.28"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.44%  "

# Row 50
$ws.Range("D50").Value = "2.289.60"
$ws.Range("E50").Value = "  +2.65%  "

# Row 51
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +76.21%  "
